$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 76212.31
$ws.Range("I86").Value = 110081.18
$ws.Range("J86").Value = 1700.8
$ws.Range("K86").Value = 110081.18
$ws.Range("L86").Value = 1700.8
$ws.Range("M86").Value = -108958.18
$ws.Range("N86").Value = -3946.8
$ws.Range("H89").Value = 76212.31
$ws.Range("I89").Value = 110081.18
$ws.Range("J89").Value = 1700.8
$ws.Range("K89").Value = 550405.8999999999
$ws.Range("L89").Value = 8504
$ws.Range("M89").Value = -544789.8999999999
$ws.Range("N89").Value = -19736
$ws.Range("H94").Value = 3500
$ws.Range("I94").Value = 3500
$ws.Range("K94").Value = 3500
$ws.Range("M94").Value = -3049
$ws.Range("H110").Value = 33451.223
$ws.Range("J110").Value = 33451.223
$ws.Range("L110").Value = 33451.223
$ws.Range("N110").Value = -41631.223
$ws.Range("H123").Value = 50780
$ws.Range("J123").Value = 50780
$ws.Range("L123").Value = 50780
$ws.Range("N123").Value = -60580
$ws.Range("H137").Value = 1864.7576
$ws.Range("I137").Value = 3165.2307
$ws.Range("J137").Value = 1019.45
$ws.Range("K137").Value = 9495.6921
$ws.Range("L137").Value = 3058.35
$ws.Range("M137").Value = -6945.6921
$ws.Range("N137").Value = -8158.35
$ws.Range("H138").Value = 2183.6191
$ws.Range("I138").Value = 1994.25
$ws.Range("J138").Value = 2271.6978
$ws.Range("K138").Value = 5982.75
$ws.Range("L138").Value = 6815.0934
$ws.Range("M138").Value = -842.75
$ws.Range("N138").Value = -17095.0934
$ws.Range("H139").Value = 49900
$ws.Range("J139").Value = 49900
$ws.Range("L139").Value = 49900
$ws.Range("N139").Value = -60180
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9585.1
$ws.Range("I32").Value = 8483.377
$ws.Range("K32").Value = 8483.377
$ws.Range("M32").Value = -8196.377
$ws.Range("H62").Value = 27000
$ws.Range("J62").Value = 27000
$ws.Range("L62").Value = 27000
$ws.Range("N62").Value = -28248
$ws.Range("H65").Value = 27000
$ws.Range("J65").Value = 27000
$ws.Range("L65").Value = 81000
$ws.Range("N65").Value = -87240
$ws.Range("H97").Value = 7457264.5
$ws.Range("I97").Value = 1112164.4
$ws.Range("J97").Value = 31251390
$ws.Range("K97").Value = 1112164.4
$ws.Range("L97").Value = 31251390
$ws.Range("M97").Value = -1111668.4
$ws.Range("N97").Value = -31252382
$ws.Range("H113").Value = 35000
$ws.Range("J113").Value = 35000
$ws.Range("L113").Value = 35000
$ws.Range("N113").Value = -43678
$ws.Range("H132").Value = 938638.4399999999
$ws.Range("I132").Value = 1987597.4
$ws.Range("J132").Value = 6230.5186
$ws.Range("K132").Value = 5962792.199999999
$ws.Range("L132").Value = 18691.5558
$ws.Range("M132").Value = -5960262.199999999
$ws.Range("N132").Value = -23751.5558
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H25").Value = 1225.4
$ws.Range("I25").Value = 1225.4
$ws.Range("K25").Value = 1225.4
$ws.Range("M25").Value = -990.4000000000001
$ws.Range("H54").Value = 1970.5
$ws.Range("I54").Value = 1970.5
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 1970.5
$ws.Range("L54").Value = 0
$ws.Range("M54").ClearContents()
$ws.Range("N54").Value = -1486.5
$ws.Range("H86").Value = 2276.375
$ws.Range("I86").Value = 1402.75
$ws.Range("J86").Value = 3150
$ws.Range("K86").Value = 1402.75
$ws.Range("L86").Value = 3150
$ws.Range("M86").Value = -279.75
$ws.Range("N86").Value = -5396
$ws.Range("H89").Value = 2276.375
$ws.Range("I89").Value = 1402.75
$ws.Range("J89").Value = 3150
$ws.Range("K89").Value = 7013.75
$ws.Range("L89").Value = 15750
$ws.Range("M89").Value = -1397.75
$ws.Range("N89").Value = -26982
$ws.Range("H94").Value = 2084605.8
$ws.Range("I94").Value = 3334482
$ws.Range("J94").Value = 1478.3334
$ws.Range("K94").Value = 3334482
$ws.Range("L94").Value = 1478.3334
$ws.Range("M94").Value = -3334031
$ws.Range("N94").Value = -2380.3334
$ws.Range("H99").Value = 2081.6
$ws.Range("I99").Value = 1100
$ws.Range("J99").Value = 2736
$ws.Range("K99").Value = 1100
$ws.Range("L99").Value = 2736
$ws.Range("M99").Value = 398
$ws.Range("N99").Value = -5732
$ws.Range("H105").Value = 1991.1818
$ws.Range("I105").Value = 1426.3572
$ws.Range("J105").Value = 2979.625
$ws.Range("K105").Value = 1426.3572
$ws.Range("L105").Value = 2979.625
$ws.Range("M105").Value = 320.6428000000001
$ws.Range("N105").Value = -6473.625
$ws.Range("H134").Value = 6794
$ws.Range("I134").Value = 2759.8667
$ws.Range("J134").Value = 10828.134
$ws.Range("K134").Value = 8279.6001
$ws.Range("L134").Value = 32484.402
$ws.Range("M134").Value = -5744.6001
$ws.Range("N134").Value = -37554.402
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H63").Value = 68271
$ws.Range("J63").Value = 68271
$ws.Range("L63").Value = 68271
$ws.Range("N63").Value = -69643
$ws.Range("H64").Value = 62271
$ws.Range("J64").Value = 62271
$ws.Range("L64").Value = 62271
$ws.Range("N64").Value = -62767
$ws.Range("H66").Value = 68271
$ws.Range("J66").Value = 68271
$ws.Range("L66").Value = 204813
$ws.Range("N66").Value = -211677
$ws.Range("H67").Value = 62271
$ws.Range("J67").Value = 62271
$ws.Range("L67").Value = 62271
$ws.Range("N67").Value = -63987
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 2214.7144
$ws.Range("I80").Value = 1500
$ws.Range("J80").Value = 4001.5
$ws.Range("K80").Value = 4500
$ws.Range("L80").Value = 12004.5
$ws.Range("M80").Value = -3564
$ws.Range("N80").Value = -13876.5
$ws.Range("H83").Value = 2214.7144
$ws.Range("I83").Value = 1500
$ws.Range("J83").Value = 4001.5
$ws.Range("K83").Value = 13500
$ws.Range("L83").Value = 36013.5
$ws.Range("M83").Value = -8820
$ws.Range("N83").Value = -45373.5
$ws.Range("H105").Value = 10001.143
$ws.Range("J105").Value = 10830.333
$ws.Range("L105").Value = 32490.999
$ws.Range("N105").Value = -37732.999
$ws.Range("H131").Value = 1060.2954
$ws.Range("I131").Value = 605
$ws.Range("J131").Value = 1132.1842
$ws.Range("K131").Value = 1815
$ws.Range("L131").Value = 3396.5526
$ws.Range("M131").Value = 3225
$ws.Range("N131").Value = -13476.5526
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 7695
$ws.Range("J5").Value = 7695
$ws.Range("L5").Value = 7695
$ws.Range("N5").Value = -7919
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 455110.12
$ws.Range("I46").Value = 467
$ws.Range("J46").Value = 769863.0600000001
$ws.Range("K46").Value = 467
$ws.Range("L46").Value = 769863.0600000001
$ws.Range("M46").Value = -279
$ws.Range("N46").Value = -770239.0600000001
$ws.Range("H64").Value = 20112.5
$ws.Range("J64").Value = 20112.5
$ws.Range("L64").Value = 20112.5
$ws.Range("N64").Value = -20562.5
$ws.Range("H67").Value = 20112.5
$ws.Range("J67").Value = 20112.5
$ws.Range("L67").Value = 20112.5
$ws.Range("N67").Value = -21672.5
$ws.Range("H82").Value = 1870.6471
$ws.Range("I82").Value = 1060.4
$ws.Range("J82").Value = 2208.25
$ws.Range("K82").Value = 1060.4
$ws.Range("L82").Value = 2208.25
$ws.Range("M82").Value = -699.4000000000001
$ws.Range("N82").Value = -2930.25
$ws.Range("H85").Value = 1870.6471
$ws.Range("I85").Value = 1060.4
$ws.Range("J85").Value = 2208.25
$ws.Range("K85").Value = 1060.4
$ws.Range("L85").Value = 2208.25
$ws.Range("M85").Value = 187.5999999999999
$ws.Range("N85").Value = -4704.25
$ws.Range("H100").Value = 3352.5
$ws.Range("I100").Value = 3117.6667
$ws.Range("K100").Value = 3117.6667
$ws.Range("M100").Value = -2576.6667
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 20349.6
$ws.Range("J63").Value = 20349.6
$ws.Range("L63").Value = 20349.6
$ws.Range("N63").Value = -21597.6
$ws.Range("H66").Value = 20349.6
$ws.Range("J66").Value = 20349.6
$ws.Range("L66").Value = 61048.8
$ws.Range("N66").Value = -67288.79999999999
$ws.Range("H96").Value = 1438.5555
$ws.Range("I96").Value = 868.6
$ws.Range("J96").Value = 2151
$ws.Range("K96").Value = 868.6
$ws.Range("L96").Value = 2151
$ws.Range("M96").Value = 504.4
$ws.Range("N96").Value = -4897
$ws.Range("H132").Value = 1435.9517
$ws.Range("I132").Value = 1289.1515
$ws.Range("J132").Value = 1603
$ws.Range("K132").Value = 3867.4545
$ws.Range("L132").Value = 4809
$ws.Range("M132").Value = -1337.4545
$ws.Range("N132").Value = -9869
